$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$times = @(
    "2021-10-05 10:50:45.582395",
    "2021-10-05 10:50:45.582408",
    "2021-10-05 10:50:45.582412",
    "2021-10-05 10:50:45.582415",
    "2021-10-05 10:50:45.582418",
    "2021-10-05 10:50:45.582421",
    "2021-10-05 10:50:45.582424",
    "2021-10-05 10:50:45.582427",
    "2021-10-05 10:50:45.582431",
    "2021-10-05 10:50:45.582434",
    "2021-10-05 10:50:45.582437",
    "2021-10-05 10:50:45.582440",
    "2021-10-05 10:50:45.582443",
    "2021-10-05 10:50:45.582445",
    "2021-10-05 10:50:45.582448",
    "2021-10-05 10:50:45.582451",
    "2021-10-05 10:50:45.582455",
    "2021-10-05 10:50:45.582458",
    "2021-10-05 10:50:45.582461",
    "2021-10-05 10:50:45.582464",
    "2021-10-05 10:50:45.582467",
    "2021-10-05 10:50:45.582470",
    "2021-10-05 10:50:45.582473",
    "2021-10-05 10:50:45.582476",
    "2021-10-05 10:50:45.582479",
    "2021-10-05 10:50:45.582482",
    "2021-10-05 10:50:45.582485",
    "2021-10-05 10:50:45.582488",
    "2021-10-05 10:50:45.582491",
    "2021-10-05 10:50:45.582494",
    "2021-10-05 10:50:45.582497",
    "2021-10-05 10:50:45.582500",
    "2021-10-05 10:50:45.582503",
    "2021-10-05 10:50:45.582506",
    "2021-10-05 10:50:45.582510",
    "2021-10-05 10:50:45.582513",
    "2021-10-05 10:50:45.582516",
    "2021-10-05 10:50:45.582519",
    "2021-10-05 10:50:45.582522",
    "2021-10-05 10:50:45.582525",
    "2021-10-05 10:50:45.582528",
    "2021-10-05 10:50:45.582532",
    "2021-10-05 10:50:45.582535",
    "2021-10-05 10:50:45.582538",
    "2021-10-05 10:50:45.582541",
    "2021-10-05 10:50:45.582543",
    "2021-10-05 10:50:45.582546",
    "2021-10-05 10:50:45.582549",
    "2021-10-05 10:50:45.582552",
    "2021-10-05 10:50:45.582555",
    "2021-10-05 10:50:45.582558"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
